$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: alexandro / bico / 1800 / 0
$ws.Range("B3").Value = "bico"
$ws.Range("C3").Value = 1800

# Update row 4: alexandro / bolinho no farol / 500 / 0 (was blank row w/ stray E4 text)
$ws.Range("A4").Value = "alexandro"
$ws.Range("B4").Value = "bolinho no farol"
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = ""

# Update row 5: alexandro / jogo playstation / 0 / 450
$ws.Range("B5").Value = "jogo playstation"
$ws.Range("D5").Value = 450

# Delete old rows 6-10 (no longer needed)
$ws.Range("A6:E10").EntireRow.Delete()
